# Appends two new paragraphs to the end of the document, after the final
# paragraph ("Club member can only be associated with one club at a time"):
#   1. A plain paragraph containing "Check queries: "
#   2. A numbered-list paragraph (same list as the rest of the document,
#      numId 1, ilvl 0) containing "9"
#
# We build each paragraph by inserting a minimal WordProcessingML package via
# Range.InsertXML on a freshly created (empty) paragraph - this lets us
# control the exact pPr/rPr contents (and avoid unwanted inheritance of the
# preceding paragraph's list formatting) instead of relying on Word's normal
# "new paragraph copies previous paragraph's formatting" behavior.

$d = $word.ActiveDocument

function New-ParagraphFromXml($afterRange, [string]$innerParagraphXml) {
    $afterRange.Collapse(0)
    $afterRange.InsertParagraphAfter()

    $target = $d.Paragraphs.Last.Range

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
             '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' +
               '<pkg:xmlData>' +
                 '<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
                   '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
                 '</Relationships>' +
               '</pkg:xmlData>' +
             '</pkg:part>' +
             '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
                 '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:body>' +
                     $innerParagraphXml +
                     '<w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr>' +
                   '</w:body>' +
                 '</w:document>' +
               '</pkg:xmlData>' +
             '</pkg:part>' +
           '</pkg:package>'

    $target.InsertXML($pkg)

    return $d.Paragraphs.Last.Range
}

# Start from the end of the current final paragraph.
$endRange = $d.Paragraphs.Last.Range

# 1) Plain paragraph: "Check queries: "
$checkQueriesXml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
                      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Check queries: </w:t></w:r>' +
                    '</w:p>'
$endRange = New-ParagraphFromXml $endRange $checkQueriesXml

# 2) Numbered list paragraph (same list used throughout the doc): "9"
$nineXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/>' +
             '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
             '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
             '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>9</w:t></w:r>' +
           '</w:p>'
$endRange = New-ParagraphFromXml $endRange $nineXml
